$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Zweite Zeile"
$ws.Range("A3").Select()
